# vpc_list.xlsx update — cn190115 / "Improve Translation for console"
#
# Content changes applied:
#  1. Corrected the English error string in C8: "subnet" -> "Subnets" (capitalised, pluralised).
#  2. Highlighted that corrected message in red and removed its word-wrap
#     (it now reads fine on one/two lines instead of three).
#  3. Added a same-styled (currently blank) column D next to it, and widened
#     the column to leave room for a future note/comment.
#  4. Shrunk row 8 back down now that the text isn't wrapping across as many
#     lines any more.
#  5. Moved the saved selection to the new D column.
#  6. Set an explicit (A4 / portrait) page setup for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ch")

$cell = $ws.Range("C8")

# 2. Re-colour / un-wrap the cell FIRST, while it still holds its original
#    text -- this keeps the cell's existing "quote-prefix" text formatting
#    (changing .Value afterwards would otherwise reset it).
$cell.WrapText = $false
$cell.Font.Color = 255

# 1. Fix the English translation text. Stage it on a scratch cell and paste
#    only the VALUE back onto C8, so the formatting we just set (and the
#    cell's quote-prefix flag) is left completely alone.
$scratch = $ws.Range("Z1")
$scratch.Value = "The VPC cannot be deleted. Please delete the Subnets of the VPC first'"
$scratch.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
$scratch.Clear()

# 3. Extend the table with a new (blank) column D carrying the same
#    formatting as C8, and give the column a sensible width.
$cell.Copy()
$ws.Range("D8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Columns.Item(4).ColumnWidth = 56

# 4. The row no longer needs to be as tall.
$ws.Rows.Item(8).RowHeight = 30

# 5. Update the remembered selection to the new column.
$ws.Range("D19").Select() | Out-Null

# 6. Explicit print setup.
$ps = $ws.PageSetup
$ps.PaperSize = 9       # xlPaperA4
$ps.Orientation = 1     # xlPortrait
